# Slide 3 ("Test Cases"), Content Placeholder 2:
#   - paragraph 4: "Make sure first card is Ace of Spades"
#                -> "Make sure first card is Ace of Hearts"
#   - paragraph 5: "Make sure last card is "
#                -> "Make sure last card " + "is King of Diamonds" (two runs)

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 4 - swap the card name in place (keeps existing run/formatting).
$para4 = $tr.Paragraphs(4)
$run4  = $para4.Runs(1)
$run4.Text = "Make sure first card is Ace of Hearts"

# Paragraph 5 - shrink the existing run to "Make sure last card " and
# append the remainder ("is King of Diamonds") as a new trailing run.
$para5 = $tr.Paragraphs(5)
$run5  = $para5.Runs(1)
$run5.Text = "Make sure last card "
[void]$para5.InsertAfter("is King of Diamonds")
